# Update vaccination data for corona report (Stand 11.5. -> Stand 12.5.)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Geimpfte Personen": headline vaccination figures
# ---------------------------------------------------------------------------
$wsGeimpft = $wb.Worksheets.Item("Geimpfte Personen")

# Header: "Stand 11.5." -> "Stand 12.5."
$wsGeimpft.Range("C1").Value = "Stand 12.5."

# Row 3 "Gesamt"
$wsGeimpft.Range("B3").Value = "24546919 (29,5 %)"
$wsGeimpft.Range("C3").Value = "28544440 (34,3 %)"
$wsGeimpft.Range("D3").Value = "4,8 PP"

# Row 4 "Nicht vollst. geimpft"
$wsGeimpft.Range("B4").Value = "17615335 (21,2 %)"
$wsGeimpft.Range("C4").Value = "20223760 (24,3 %)"
$wsGeimpft.Range("D4").Value = "3,1 PP"

# Row 5 "Vollst. geimpft"
$wsGeimpft.Range("B5").Value = "6931584 ( 8,3 %)"
$wsGeimpft.Range("C5").Value = "8320680 (10,0 %)"
$wsGeimpft.Range("D5").Value = "1,7 PP"

# ---------------------------------------------------------------------------
# Sheet "Regional Geimpfte": per-region "Gesamt min. 1x" / "Gesamt vollst." %
# ---------------------------------------------------------------------------
$wsRegional = $wb.Worksheets.Item("Regional Geimpfte")

# row -> (Gesamt min. 1x %, Gesamt vollst. %)
$wsRegional.Range("C2").Value  = "34,3"
$wsRegional.Range("D2").Value  = "10,0"

$wsRegional.Range("C3").Value  = "33,7"
$wsRegional.Range("D3").Value  = " 9,7"

$wsRegional.Range("C4").Value  = "34,9"
$wsRegional.Range("D4").Value  = " 9,5"

$wsRegional.Range("C5").Value  = "31,3"
$wsRegional.Range("D5").Value  = "11,6"

$wsRegional.Range("C6").Value  = "30,1"
$wsRegional.Range("D6").Value  = "10,6"

$wsRegional.Range("C7").Value  = "35,7"
$wsRegional.Range("D7").Value  = "12,9"

$wsRegional.Range("C8").Value  = "34,0"
$wsRegional.Range("D8").Value  = " 9,1"

$wsRegional.Range("C9").Value  = "34,1"
$wsRegional.Range("D9").Value  = " 9,0"

$wsRegional.Range("C10").Value = "36,6"
$wsRegional.Range("D10").Value = " 9,3"

$wsRegional.Range("C11").Value = "35,6"
$wsRegional.Range("D11").Value = " 8,6"

$wsRegional.Range("C12").Value = "36,3"
$wsRegional.Range("D12").Value = " 9,1"

$wsRegional.Range("C13").Value = "32,9"
$wsRegional.Range("D13").Value = "11,7"

$wsRegional.Range("C14").Value = "38,7"
$wsRegional.Range("D14").Value = "10,2"

$wsRegional.Range("C15").Value = "29,8"
$wsRegional.Range("D15").Value = "13,7"

$wsRegional.Range("C16").Value = "33,6"
$wsRegional.Range("D16").Value = "10,0"

$wsRegional.Range("C17").Value = "31,9"
$wsRegional.Range("D17").Value = "11,9"

$wsRegional.Range("C18").Value = "31,2"
$wsRegional.Range("D18").Value = "14,4"

# ---------------------------------------------------------------------------
# Sheet "Impfstoffdosen": doses delivered/used per vaccine manufacturer
# ---------------------------------------------------------------------------
$wsDosen = $wb.Worksheets.Item("Impfstoffdosen")

# These cells hold plain digit strings (e.g. "17014901") that must stay text
# (as in the source file) instead of being auto-converted to numbers by
# Excel. Mark them as Text-formatted first so the value assignment below
# keeps them as strings.
$textCells = @("B3","C3","B4","C4","B7","C7","B8","C8","B11","C11","B12","C12")
foreach ($addr in $textCells) {
    $wsDosen.Range($addr).NumberFormat = "@"
}

# Biontech/Pfizer
$wsDosen.Range("B2").Value  = "23399097 (74,4 %)"
$wsDosen.Range("C2").Value  = "27168667 (73,8 %)"
$wsDosen.Range("B3").Value  = "17014901"
$wsDosen.Range("C3").Value  = "19701594"
$wsDosen.Range("B4").Value  = "6384196"
$wsDosen.Range("C4").Value  = "7467073"

# Moderna
$wsDosen.Range("B6").Value  = "1932692 ( 6,1 %)"
$wsDosen.Range("C6").Value  = "2463861 ( 6,7 %)"
$wsDosen.Range("B7").Value  = "1482621"
$wsDosen.Range("C7").Value  = "1870450"
$wsDosen.Range("B8").Value  = "450071"
$wsDosen.Range("C8").Value  = "593411"

# AstraZeneca
$wsDosen.Range("B10").Value = "6127132 (19,5 %)"
$wsDosen.Range("C10").Value = "7176720 (19,5 %)"
$wsDosen.Range("B11").Value = "6039606"
$wsDosen.Range("C11").Value = "6944460"
$wsDosen.Range("B12").Value = "87526"
$wsDosen.Range("C12").Value = "232260"

# Johnson&Johnson
$wsDosen.Range("B14").Value = "9791 ( 0,0 %)"
$wsDosen.Range("C14").Value = "27936 ( 0,1 %)"
